$wb = $excel.ActiveWorkbook

# Data refresh for 2024-12-02: update 2024 year-to-date (column K) crime
# counts (and a couple of adjacent 2023 column-J corrections) across the
# citywide, by-neighborhood, and per-neighborhood detail sheets.

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7373
$ws.Range("K3").Value = 7643
$ws.Range("J4").Value = 1736
$ws.Range("K4").Value = 1608
$ws.Range("K5").Value = 539
$ws.Range("K6").Value = 8459
$ws.Range("J7").Value = 26959
$ws.Range("K7").Value = 25622

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 503
$ws.Range("K4").Value = 96
$ws.Range("K6").Value = 557
$ws.Range("K7").Value = 1665

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 188
$ws.Range("K6").Value = 126
$ws.Range("K7").Value = 541

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 386
$ws.Range("K6").Value = 343
$ws.Range("K7").Value = 1081

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 282
$ws.Range("K6").Value = 257
$ws.Range("K7").Value = 854

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 187
$ws.Range("K7").Value = 775
$ws.Range("K8").Value = 1665
$ws.Range("K10").Value = 153
$ws.Range("K11").Value = 462
$ws.Range("K15").Value = 262
$ws.Range("K16").Value = 60
$ws.Range("K18").Value = 170
$ws.Range("K19").Value = 740
$ws.Range("K24").Value = 82
$ws.Range("K25").Value = 118
$ws.Range("K29").Value = 1409
$ws.Range("K31").Value = 298
$ws.Range("K33").Value = 1081
$ws.Range("K34").Value = 147
$ws.Range("K37").Value = 854
$ws.Range("K43").Value = 210
$ws.Range("K48").Value = 325
$ws.Range("K49").Value = 142
$ws.Range("K52").Value = 668
$ws.Range("K54").Value = 500
$ws.Range("K55").Value = 279
$ws.Range("K57").Value = 102
$ws.Range("K63").Value = 72
$ws.Range("K64").Value = 152
$ws.Range("K66").Value = 77
$ws.Range("K67").Value = 1001
$ws.Range("K78").Value = 315
$ws.Range("K80").Value = 96
$ws.Range("K83").Value = 541
$ws.Range("K84").Value = 203
$ws.Range("K85").Value = 1172
$ws.Range("K89").Value = 386
$ws.Range("K90").Value = 250
$ws.Range("K91").Value = 305
$ws.Range("K92").Value = 94
$ws.Range("K94").Value = 340
$ws.Range("J96").Value = 293
$ws.Range("K98").Value = 135
$ws.Range("J101").Value = 26959
$ws.Range("K101").Value = 25622

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 91
$ws.Range("K7").Value = 298

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 362
$ws.Range("K4").Value = 58
$ws.Range("K7").Value = 1001

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 203

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 114
$ws.Range("K4").Value = 33
$ws.Range("K6").Value = 274
$ws.Range("K7").Value = 500

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 400
$ws.Range("K3").Value = 502
$ws.Range("K6").Value = 409
$ws.Range("K7").Value = 1409

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 80
$ws.Range("K7").Value = 325

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 245
$ws.Range("K7").Value = 740

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K3").Value = 48
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K4").Value = 29
$ws.Range("K7").Value = 315

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 80
$ws.Range("K6").Value = 101
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 293

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 78
$ws.Range("K7").Value = 305

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 44
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 255
$ws.Range("K7").Value = 775

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 57
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 160
$ws.Range("K7").Value = 340

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K2").Value = 46
$ws.Range("K7").Value = 118

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 101
$ws.Range("K7").Value = 262

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K2").Value = 26
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 117
$ws.Range("K7").Value = 462

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 94

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 119
$ws.Range("K7").Value = 386

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 250

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K4").Value = 28
$ws.Range("K7").Value = 210

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 385
$ws.Range("K5").Value = 33
$ws.Range("K6").Value = 292
$ws.Range("K7").Value = 1172

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 185
$ws.Range("K7").Value = 668

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 60
